# Generate Report for Handback
# Update the handoff/handback timestamps recorded for the zh-cn and de-de
# language sheets to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 2 and 3 share the same handoff/handback pair of values
$wsZhCn.Range("E2:E3").Value = "2016-03-13 08:20:39"
$wsZhCn.Range("H2:H3").Value = "2016-03-13 08:20:59"

# de-de sheet: rows 2 and 3 share the same handoff/handback pair of values
$wsDeDe.Range("E2:E3").Value = "2016-03-13 08:20:43"
$wsDeDe.Range("H2:H3").Value = "2016-03-13 08:21:05"
